$d = $word.ActiveDocument

$replacements = @(
    @("731÷4=182, 3", "648÷9=72, 0"),
    @("945÷2=472, 1", "474÷7=67, 5"),
    @("702÷7=100, 2", "742÷9=82, 4"),
    @("573÷9=63, 6", "928÷8=116, 0"),
    @("904÷7=129, 1", "361÷9=40, 1"),
    @("318÷4=79, 2", "148÷2=74, 0"),
    @("774÷7=110, 4", "380÷5=76, 0"),
    @("447÷9=49, 6", "215÷5=43, 0"),
    @("653÷4=163, 1", "946÷8=118, 2"),
    @("447÷7=63, 6", "366÷2=183, 0"),
    @("822÷3=274, 0", "204÷9=22, 6"),
    @("303÷4=75, 3", "455÷9=50, 5"),
    @("835÷2=417, 1", "648÷6=108, 0"),
    @("309÷7=44, 1", "879÷9=97, 6"),
    @("162÷4=40, 2", "504÷3=168, 0"),
    @("240÷3=80, 0", "640÷2=320, 0"),
    @("527÷8=65, 7", "474÷3=158, 0"),
    @("643÷5=128, 3", "790÷3=263, 1"),
    @("350÷5=70, 0", "980÷5=196, 0"),
    @("420÷7=60, 0", "913÷7=130, 3"),
    @("530÷3=176, 2", "628÷9=69, 7"),
    @("228÷4=57, 0", "365÷8=45, 5"),
    @("929÷2=464, 1", "917÷3=305, 2"),
    @("820÷9=91, 1", "777÷5=155, 2"),
    @("136÷2=68, 0", "248÷6=41, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
